$wb = $excel.ActiveWorkbook

# --- Sheet: Simple Fields ---
$ws = $wb.Worksheets.Item("Simple Fields")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("R2").NumberFormat = "@"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("B2").Value = "0.9359569"
$ws.Range("D2").Value = "0.91714096"
$ws.Range("F2").Value = "0.9687482"
$ws.Range("H2").Value = "0.68830854"
$ws.Range("J2").Value = "0.91500014"
$ws.Range("L2").Value = "0.9930566"
$ws.Range("N2").Value = "0.95059097"
$ws.Range("P2").Value = "0.9489661"
$ws.Range("R2").Value = "0.9461934"
$ws.Range("T2").Value = "0.6617088"

# --- Sheet: Simple Fields - Formatted ---
$ws = $wb.Worksheets.Item("Simple Fields - Formatted")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("R2").NumberFormat = "@"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("B2").Value = "0.9359569"
$ws.Range("D2").Value = "0.91714096"
$ws.Range("F2").Value = "0.9687482"
$ws.Range("H2").Value = "0.68830854"
$ws.Range("J2").Value = "0.91500014"
$ws.Range("L2").Value = "0.9930566"
$ws.Range("N2").Value = "0.95059097"
$ws.Range("P2").Value = "0.9489661"
$ws.Range("R2").Value = "0.9461934"
$ws.Range("T2").Value = "0.6617088"

# --- Sheet: Items ---
$ws = $wb.Worksheets.Item("Items")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("A2").Value = "green onion Pancakes NVIW@If (1)"
$ws.Range("B2").Value = "0.6617088"
$ws.Range("D2").Value = "0.9654466"
$ws.Range("H2").Value = "0.96437496"
$ws.Range("A3").Value = "Pan Fried Leek Dumplings #j (2)"
$ws.Range("B3").Value = "0.72182536"
$ws.Range("D3").Value = "0.92594826"
$ws.Range("H3").Value = "0.9631321"
$ws.Range("A4").Value = "Pork Xiao Long Bao(10) MJINAME(10)"
$ws.Range("B4").Value = "0.74614596"
$ws.Range("D4").Value = "0.8857973"
$ws.Range("H4").Value = "0.97184485"
$ws.Range("A5").Value = "Q-BAO (5) HENE] (5)"
$ws.Range("B5").Value = "0.8768804"
$ws.Range("D5").Value = "0.8922952"
$ws.Range("H5").Value = "0.96625626"
$ws.Range("A6").Value = "Chicken potstickers #ERJWAUA(6)"
$ws.Range("B6").Value = "0.9518383"
$ws.Range("D6").Value = "0.92076415"
$ws.Range("H6").Value = "0.9781219"
$ws.Range("A7").Value = "Tomato Mushroom Steamed dumpli ¿ (6)"
$ws.Range("B7").Value = "0.929742"
$ws.Range("D7").Value = "0.93721306"
$ws.Range("H7").Value = "0.96639776"
$ws.Range("A8").Value = "Zucchini shrimp dumplings jJUUANUIC"
$ws.Range("B8").Value = "0.97735536"
$ws.Range("D8").Value = "0.9401775"
$ws.Range("H8").Value = "0.9734021"
$ws.Range("A9").Value = "beef stew nodle soup (Non Spicy 84pJ(T#)"
$ws.Range("B9").Value = "0.8927655"
$ws.Range("D9").Value = "0.91823614"
$ws.Range("H9").Value = "0.9591488"
$ws.Range("A10").Value = "dandan noodle ttIÉÍ"
$ws.Range("B10").Value = "0.97385263"
$ws.Range("D10").Value = "0.909689"
$ws.Range("H10").Value = "0.9875826"
$ws.Range("A11").Value = "banana naan bread EATA"
$ws.Range("B11").Value = "0.9784949"
$ws.Range("D11").Value = "0.91455704"
$ws.Range("H11").Value = "0.9884863"
$ws.Range("A12").Value = "house made plum juice piumit"
$ws.Range("B12").Value = "0.9537653"
$ws.Range("D12").Value = "0.90424895"
$ws.Range("H12").Value = "0.9838298"

# --- Sheet: Items - Formatted ---
$ws = $wb.Worksheets.Item("Items - Formatted")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("A2").Value = "green onion Pancakes NVIW@If (1)"
$ws.Range("B2").Value = "0.6617088"
$ws.Range("D2").Value = "0.9654466"
$ws.Range("H2").Value = "0.96437496"
$ws.Range("A3").Value = "Pan Fried Leek Dumplings #j (2)"
$ws.Range("B3").Value = "0.72182536"
$ws.Range("D3").Value = "0.92594826"
$ws.Range("H3").Value = "0.9631321"
$ws.Range("A4").Value = "Pork Xiao Long Bao(10) MJINAME(10)"
$ws.Range("B4").Value = "0.74614596"
$ws.Range("D4").Value = "0.8857973"
$ws.Range("H4").Value = "0.97184485"
$ws.Range("A5").Value = "Q-BAO (5) HENE] (5)"
$ws.Range("B5").Value = "0.8768804"
$ws.Range("D5").Value = "0.8922952"
$ws.Range("H5").Value = "0.96625626"
$ws.Range("A6").Value = "Chicken potstickers #ERJWAUA(6)"
$ws.Range("B6").Value = "0.9518383"
$ws.Range("D6").Value = "0.92076415"
$ws.Range("H6").Value = "0.9781219"
$ws.Range("A7").Value = "Tomato Mushroom Steamed dumpli ¿ (6)"
$ws.Range("B7").Value = "0.929742"
$ws.Range("D7").Value = "0.93721306"
$ws.Range("H7").Value = "0.96639776"
$ws.Range("A8").Value = "Zucchini shrimp dumplings jJUUANUIC"
$ws.Range("B8").Value = "0.97735536"
$ws.Range("D8").Value = "0.9401775"
$ws.Range("H8").Value = "0.9734021"
$ws.Range("A9").Value = "beef stew nodle soup (Non Spicy 84pJ(T#)"
$ws.Range("B9").Value = "0.8927655"
$ws.Range("D9").Value = "0.91823614"
$ws.Range("H9").Value = "0.9591488"
$ws.Range("A10").Value = "dandan noodle ttIÉÍ"
$ws.Range("B10").Value = "0.97385263"
$ws.Range("D10").Value = "0.909689"
$ws.Range("H10").Value = "0.9875826"
$ws.Range("A11").Value = "banana naan bread EATA"
$ws.Range("B11").Value = "0.9784949"
$ws.Range("D11").Value = "0.91455704"
$ws.Range("H11").Value = "0.9884863"
$ws.Range("A12").Value = "house made plum juice piumit"
$ws.Range("B12").Value = "0.9537653"
$ws.Range("D12").Value = "0.90424895"
$ws.Range("H12").Value = "0.9838298"
